# Apply targeted cell value updates to Sheet1 based on the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = 2.3
$ws.Range("I2").Value = 2.7
$ws.Range("N2").Value = 3
$ws.Range("P2").Value = 1.66
$ws.Range("U2").Value = 1.95
$ws.Range("V2").Value = 1.59
$ws.Range("J3").Value = 1.2
$ws.Range("U3").Value = 1.83
$ws.Range("N4").Value = 3.45
$ws.Range("O4").Value = 1.23
$ws.Range("P4").Value = 2.12
$ws.Range("Q4").Value = 1.72
$ws.Range("H6").Value = 2.48
$ws.Range("I6").Value = 2.74
$ws.Range("J6").Value = 3.5
$ws.Range("S6").Value = 2.84
$ws.Range("T6").Value = 1.7
$ws.Range("V6").Value = 1.57
$ws.Range("AA6").Value = 980
$ws.Range("AB6").Value = 16.5
$ws.Range("AE6").Value = 30
$ws.Range("AF6").Value = 22
$ws.Range("AI6").Value = 38
$ws.Range("AK6").Value = 32
$ws.Range("AN6").Value = 23
$ws.Range("F7").Value = 4
$ws.Range("F8").Value = 2.68
$ws.Range("H8").Value = 2.38
$ws.Range("N8").Value = 2.78
$ws.Range("Q8").Value = 2
$ws.Range("F9").Value = 1.84
$ws.Range("G9").Value = 1.94
$ws.Range("I9").Value = 5.1
$ws.Range("K9").Value = 4
$ws.Range("W9").Value = 2.06
$ws.Range("F10").Value = 2.42
$ws.Range("G10").Value = 2.6
$ws.Range("H10").Value = 2.8
$ws.Range("I10").Value = 3
$ws.Range("K10").Value = 4
$ws.Range("L10").Value = 1.28
$ws.Range("O10").Value = 1.24
$ws.Range("P10").Value = 2.2
$ws.Range("Q10").Value = 1.72
$ws.Range("S10").Value = 2.74
$ws.Range("V10").Value = 1.5
$ws.Range("W10").Value = 1.63
$ws.Range("X10").Value = 20
$ws.Range("Y10").Value = 15.5
$ws.Range("Z10").Value = 23
$ws.Range("AB10").Value = 14
$ws.Range("AF10").Value = 19
$ws.Range("F11").Value = 4.2
$ws.Range("H11").Value = 1.79
$ws.Range("K11").Value = 4.7
$ws.Range("T11").Value = 1.55
$ws.Range("X11").Value = 34
$ws.Range("AB11").Value = 29
$ws.Range("AC11").Value = 12.5
$ws.Range("AD11").Value = 11
$ws.Range("AH11").Value = 19
$ws.Range("F12").Value = 2.8
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = 2.58
$ws.Range("I12").Value = 2.78
$ws.Range("L12").Value = 1.42
$ws.Range("P12").Value = 1.85
$ws.Range("U12").Value = 2.1
$ws.Range("V12").Value = 1.56
$ws.Range("W12").Value = 1.5
$ws.Range("Y12").Value = 11.5
$ws.Range("Z12").Value = 21
$ws.Range("AC12").Value = 8
$ws.Range("AJ12").Value = 50
$ws.Range("AN12").Value = 36
$ws.Range("AO12").Value = 27
$ws.Range("I13").Value = 5.2
$ws.Range("J13").Value = 3.9
$ws.Range("K13").Value = 5
$ws.Range("L13").Value = 1.01
$ws.Range("T13").Value = 1.65
$ws.Range("U13").Value = 2.08
$ws.Range("V13").Value = 1.25
$ws.Range("AH13").Value = 18.5
$ws.Range("AN13").Value = 10.5
$ws.Range("F14").Value = 1.91
$ws.Range("K14").Value = 4.1
$ws.Range("T14").Value = 1.94
